$d = $word.ActiveDocument

$replacements = @(
    @("47×71=3337", "88×27=2376"),
    @("88×75=6600", "79×61=4819"),
    @("96×23=2208", "42×18=756"),
    @("36×71=2556", "75×81=6075"),
    @("30×66=1980", "27×15=405"),
    @("50×65=3250", "41×38=1558"),
    @("58×36=2088", "72×24=1728"),
    @("44×73=3212", "16×82=1312"),
    @("47×59=2773", "17×19=323"),
    @("92×53=4876", "47×68=3196"),
    @("68×19=1292", "63×50=3150"),
    @("91×26=2366", "97×64=6208"),
    @("77×32=2464", "11×67=737"),
    @("44×89=3916", "95×34=3230"),
    @("19×26=494", "77×21=1617"),
    @("68×75=5100", "61×85=5185"),
    @("17×34=578", "25×66=1650"),
    @("41×94=3854", "19×99=1881"),
    @("85×62=5270", "68×62=4216"),
    @("97×27=2619", "69×65=4485"),
    @("86×78=6708", "76×12=912"),
    @("11×82=902", "42×54=2268"),
    @("30×54=1620", "24×82=1968"),
    @("54×96=5184", "23×67=1541"),
    @("35×79=2765", "20×16=320")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
